$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 60) to the log, reusing the formatting
# (date / day-name / lat-long number formats) from the last existing
# row (59) so the new row's styles match the rest of the table.
$ws.Range("A59:K59").Copy()
$ws.Range("A60").PasteSpecial(-4122)

# Date serial 45755 = 2025-04-08 (a Tuesday)
$ws.Range("A60").Value = 45755
$ws.Range("B60").Value = "Tuesday"
$ws.Range("C60").Value = "Derek George"
$ws.Range("D60").Value = "Male"
$ws.Range("E60").Value = "50 to 60 Years"
$ws.Range("F60").Value = "SKDP"
$ws.Range("G60").Value = "Julie Brosnan"
$ws.Range("H60").Value = "Kerry"
$ws.Range("I60").Value = 52.046317
$ws.Range("J60").Value = -9.339747
$ws.Range("K60").Value = 1
